$wb = $excel.ActiveWorkbook

# --- Sheet "BuscaLupa": add a new search entry ("laptop" / "HP Stream") ---
# at the top of the history, shifting the previous entries one column to
# the right and dropping the old "chiclete" leftover value.
$wsLupa = $wb.Worksheets.Item("BuscaLupa")

$oldA1 = $wsLupa.Range("A1").Value2
$oldB1 = $wsLupa.Range("B1").Value2
$oldA2 = $wsLupa.Range("A2").Value2

$wsLupa.Range("C1").Value = $oldB1
$wsLupa.Range("B1").Value = $oldA1
$wsLupa.Range("A1").Value = "laptop"

$wsLupa.Range("B2").Value = $oldA2
$wsLupa.Range("A2").Value = "HP Stream"

$wsLupa.Range("C1").Select() | Out-Null

# --- Sheet "Contas": update the sample password value used in row 2, and
# restore it as the active/selected sheet (as it was before the edit). ---
$wsContas = $wb.Worksheets.Item("Contas")
$wsContas.Range("A2").Value = "Miiera5"
$wsContas.Activate() | Out-Null

$wb.Save()
